{"js": "// Append two new narrative paragraphs (each preceded by a blank separator\n// paragraph) to the end of the progress-log document, describing the\n// AccommDetailsView map-tooltip feature and the BookingViewForm /\n// ListItemSelectedView refactor.\n\nconst body = context.document.body;\n\n// Locate the very last paragraph currently in the document\n// (\"...modifier\u00e9nek.\") so the new content is appended after it.\nconst lastParagraph = body.paragraphs.getLast();\n\n// Blank paragraph used as a visual separator, then the first new\n// paragraph describing the AccommDetailsView / MapKit tooltip work.\nconst blankBeforeFirst = lastParagraph.insertParagraph(\"\", \"After\");\nconst firstNewParagraph = blankBeforeFirst.insertParagraph(\n  \"Ezek mellett az AccommDetailsView k\u00e9perny\u0151re kitettem a \u201eF\u00f6ldrajzi helyzet\u201d label mell\u00e9 egy kis k\u00e9k tooltipet, amelyre kattintva a felhaszn\u00e1l\u00f3 egy alert form\u00e1j\u00e1ban taj\u00e9koztat\u00e1st kap arr\u00f3l, hogy a f\u00f6ldrajzi koordin\u00e1t\u00e1kra tappolva azok megnyit\u00e1sra ker\u00fclnek a T\u00e9rk\u00e9pek alkalmaz\u00e1sban, \u00e9s egy gombost\u0171 form\u00e1j\u00e1ban megtekintheti, hogy pontosan hol tal\u00e1lht\u00f3 az aktu\u00e1lisan megtekintett sz\u00e1ll\u00e1s a t\u00e9rk\u00e9pen. Ehhez a f\u00f6ldrajzi koordin\u00e1t\u00e1kat tartalmaz\u00f3 Text objektumon elhelyeztem egy \u201eonTapGesture\u201d modifiert, amely megh\u00edv egy \u00e1ltalam defini\u00e1lt \u201eopenLocationOnMap\u201d callback f\u00fcggv\u00e9nyt a koordin\u00e1t\u00e1kkal. Ez a f\u00fcggv\u00e9ny megfelel\u0151 t\u00edpus\u00fa \u00e9rt\u00e9keket hoz l\u00e9tre a koordin\u00e1t\u00e1kb\u00f3l, majd ezek alapj\u00e1n l\u00e9trehoz egy 10 km-es sugar\u00fa r\u00e9gi\u00f3t a meghat\u00e1rozott pont k\u00f6r\u00fcl, amelyet megjelen\u00edt a T\u00e9rk\u00e9pek alkalmaz\u00e1sban. Ehhez a funkci\u00f3hoz beimport\u00e1ltam a MapKit k\u00f6nyvt\u00e1rat a swift f\u00e1jlban.\",\n  \"After\"\n);\n\n// Blank paragraph used as a visual separator, then the second new\n// paragraph describing the BookingViewForm / ListItemSelectedView work.\nconst blankBeforeSecond = firstNewParagraph.insertParagraph(\"\", \"After\");\nconst secondNewParagraph = blankBeforeSecond.insertParagraph(\n  \"Egy kisebb refaktor\u00e1l\u00e1son esett \u00e1t a BookingViewForm k\u00e9perny\u0151: lehet\u0151v\u00e9 tettem, hogy a keres\u00e9s ind\u00edt\u00e1sakor egyszerre t\u00f6bb fizet\u00e9si m\u00f3d is kiv\u00e1laszthat\u00f3 legyen. Ehhez m\u00f3dos\u00edtani kellett az \u201eOtherNeedsListItemSelectedView\u201d oszt\u00e1lyt: ahhoz, hogy \u00e1ltal\u00e1nosabb felhaszn\u00e1l\u00e1s\u00fa legyen \u00e9s ezt a neve is t\u00fckr\u00f6zze, \u00e1tneveztem \u201eListItemSelectedView\u201d-ra \u00e9s \u00e1tmozgattam a utils mapp\u00e1ba. Ezt a v\u00e1ltoz\u00e1st le kellett k\u00f6vetni t\u00f6bb helyen is a k\u00f3dban, illetve l\u00e9trehoztam egy ListItem oszt\u00e1lyt, amely \u00e1tvette az eddigi \u201eOtherNeed\u201d oszt\u00e1ly szerep\u00e9t a lista l\u00e9trehoz\u00e1sban. Az \u201eOtherNeed\u201d oszt\u00e1lyb\u00f3l pedig elt\u00e1vol\u00edtottam az \u201eis_selected\u201dproperty-t, amelyre kiz\u00e1r\u00f3lag a lista l\u00e9trehoz\u00e1s miatt volt eddig sz\u00fcks\u00e9g. Ezen k\u00edv\u00fcl megjelenik a sz\u00e1ll\u00e1s r\u00e9szletei oldalon is, hogy milyen elfogadott fizet\u00e9si m\u00f3dok \u00e1llnak rendelkez\u00e9sre az adott sz\u00e1ll\u00e1son.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Append two new narrative paragraphs (each preceded by a blank separator\n# paragraph) to the end of the progress-log document, describing the\n# AccommDetailsView map-tooltip feature and the BookingViewForm /\n# ListItemSelectedView refactor.\n\n$d = $word.ActiveDocument\n\n# --- Blank separator paragraph after the current last paragraph\n#     (\"...modifier\u00e9nek.\") ---\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n# --- New paragraph describing the AccommDetailsView / MapKit tooltip work ---\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)\n$r.InsertAfter(\"Ezek mellett az AccommDetailsView k\u00e9perny\u0151re kitettem a \u201eF\u00f6ldrajzi helyzet\u201d label mell\u00e9 egy kis k\u00e9k tooltipet, amelyre kattintva a felhaszn\u00e1l\u00f3 egy alert form\u00e1j\u00e1ban taj\u00e9koztat\u00e1st kap arr\u00f3l, hogy a f\u00f6ldrajzi koordin\u00e1t\u00e1kra tappolva azok megnyit\u00e1sra ker\u00fclnek a T\u00e9rk\u00e9pek alkalmaz\u00e1sban, \u00e9s egy gombost\u0171 form\u00e1j\u00e1ban megtekintheti, hogy pontosan hol tal\u00e1lht\u00f3 az aktu\u00e1lisan megtekintett sz\u00e1ll\u00e1s a t\u00e9rk\u00e9pen. Ehhez a f\u00f6ldrajzi koordin\u00e1t\u00e1kat tartalmaz\u00f3 Text objektumon elhelyeztem egy \u201eonTapGesture\u201d modifiert, amely megh\u00edv egy \u00e1ltalam defini\u00e1lt \u201eopenLocationOnMap\u201d callback f\u00fcggv\u00e9nyt a koordin\u00e1t\u00e1kkal. Ez a f\u00fcggv\u00e9ny megfelel\u0151 t\u00edpus\u00fa \u00e9rt\u00e9keket hoz l\u00e9tre a koordin\u00e1t\u00e1kb\u00f3l, majd ezek alapj\u00e1n l\u00e9trehoz egy 10 km-es sugar\u00fa r\u00e9gi\u00f3t a meghat\u00e1rozott pont k\u00f6r\u00fcl, amelyet megjelen\u00edt a T\u00e9rk\u00e9pek alkalmaz\u00e1sban. Ehhez a funkci\u00f3hoz beimport\u00e1ltam a MapKit k\u00f6nyvt\u00e1rat a swift f\u00e1jlban.\")\n\n# --- Blank separator paragraph ---\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n# --- New paragraph describing the BookingViewForm / ListItemSelectedView work ---\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)\n$r.InsertAfter(\"Egy kisebb refaktor\u00e1l\u00e1son esett \u00e1t a BookingViewForm k\u00e9perny\u0151: lehet\u0151v\u00e9 tettem, hogy a keres\u00e9s ind\u00edt\u00e1sakor egyszerre t\u00f6bb fizet\u00e9si m\u00f3d is kiv\u00e1laszthat\u00f3 legyen. Ehhez m\u00f3dos\u00edtani kellett az \u201eOtherNeedsListItemSelectedView\u201d oszt\u00e1lyt: ahhoz, hogy \u00e1ltal\u00e1nosabb felhaszn\u00e1l\u00e1s\u00fa legyen \u00e9s ezt a neve is t\u00fckr\u00f6zze, \u00e1tneveztem \u201eListItemSelectedView\u201d-ra \u00e9s \u00e1tmozgattam a utils mapp\u00e1ba. Ezt a v\u00e1ltoz\u00e1st le kellett k\u00f6vetni t\u00f6bb helyen is a k\u00f3dban, illetve l\u00e9trehoztam egy ListItem oszt\u00e1lyt, amely \u00e1tvette az eddigi \u201eOtherNeed\u201d oszt\u00e1ly szerep\u00e9t a lista l\u00e9trehoz\u00e1sban. Az \u201eOtherNeed\u201d oszt\u00e1lyb\u00f3l pedig elt\u00e1vol\u00edtottam az \u201eis_selected\u201dproperty-t, amelyre kiz\u00e1r\u00f3lag a lista l\u00e9trehoz\u00e1s miatt volt eddig sz\u00fcks\u00e9g. Ezen k\u00edv\u00fcl megjelenik a sz\u00e1ll\u00e1s r\u00e9szletei oldalon is, hogy milyen elfogadott fizet\u00e9si m\u00f3dok \u00e1llnak rendelkez\u00e9sre az adott sz\u00e1ll\u00e1son.\")\n"}
